$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 943
$ws.Range("F3").Value = 1036
$ws.Range("F4").Value = 829
$ws.Range("F5").Value = 895
$ws.Range("F7").Value = 730
$ws.Range("F9").Value = 1337
$ws.Range("F10").Value = 771
$ws.Range("F11").Value = 434
$ws.Range("F12").Value = 576
$ws.Range("F13").Value = 193
$ws.Range("F14").Value = 78
$ws.Range("F15").Value = 78
$ws.Range("F16").Value = 1336
$ws.Range("F17").Value = 154
$ws.Range("F19").Value = 442
$ws.Range("F20").Value = 16
$ws.Range("F23").Value = 612
$ws.Range("F24").Value = 172
$ws.Range("F25").Value = 675
$ws.Range("F26").Value = 43
$ws.Range("F27").Value = 1185
$ws.Range("F28").Value = 26

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 351
$ws.Range("F5").Value = 647
$ws.Range("F7").Value = 261
$ws.Range("F8").Value = 60
$ws.Range("F11").Value = 123

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 351
$ws.Range("F4").Value = 943
$ws.Range("F5").Value = 1036
$ws.Range("F6").Value = 829
$ws.Range("F7").Value = 895
$ws.Range("F9").Value = 730
$ws.Range("F11").Value = 1337
$ws.Range("F12").Value = 771
$ws.Range("F15").Value = 434
$ws.Range("F16").Value = 576
$ws.Range("F17").Value = 647
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 78
$ws.Range("F20").Value = 78
$ws.Range("F21").Value = 1336
$ws.Range("F23").Value = 154
$ws.Range("F25").Value = 442
$ws.Range("F26").Value = 16
$ws.Range("F29").Value = 261
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 612
$ws.Range("F34").Value = 123
$ws.Range("F35").Value = 123
$ws.Range("F36").Value = 172
$ws.Range("F37").Value = 675
$ws.Range("F38").Value = 43
$ws.Range("F39").Value = 1185
$ws.Range("F40").Value = 26
